$d = $word.ActiveDocument

# The commit removes the "Requisitos" section (a Heading2 paragraph
# followed by a ListBullet paragraph reading
# "LOQ4073 -  Química Geral II  (Requisito fraco)") that used to sit at
# the very end of the document, right before the sectPr. Locate that
# heading paragraph and delete everything from its start through the
# end of the document content, which removes both paragraphs (and their
# paragraph marks) while leaving the Bibliografia text intact.

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Requisitos") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Requisitos' heading paragraph"
}

$rng = $d.Range($target.Range.Start, $d.Content.End)
$rng.Delete()
